$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the RF (column I) values for rows 20 through 53 to reflect the
# 2025 data / RF recalculation.
$ws.Range("I20:I53").Value = 11.95527272727273
